$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting (styles, fill, borders) of the last existing data row (25)
# down into the two new rows so the new rows look consistent with the rest
# of the table.
$ws.Range("A25:E25").Copy($ws.Range("A26:E26"))
$ws.Range("A25:E25").Copy($ws.Range("A27:E27"))

# ---- Row 26 : Notifications0025 / OPQA-4493 ----
$ws.Range("A26").Value2 = "Notifications0025"
$ws.Range("B26").Value2 = "OPQA-4493"
$ws.Range("C26").Value2 = "Verify that user navigate to record view page of the article while clicking article in trending section from Newsfeed page."
$ws.Range("D26").Value2 = "Y"
$ws.Range("E26").Value2 = "PASS"

# ---- Row 27 : Notifications0026 / OPQA-4499 ----
$ws.Range("A27").Value2 = "Notifications0026"
$ws.Range("C27").Value2 = "Verify that user navigate to record view page of the post while clicking post in trending section from Newsfeed page."
$ws.Range("B27").Value2 = "OPQA-4499"
$ws.Range("D27").Value2 = "Y"
$ws.Range("E27").Value2 = "PASS"

# Update the view so the newly added row is shown/selected, matching what a
# user would see after scrolling down and clicking on the new row.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 13
$ws.Range("C27").Select() | Out-Null
